# issue #5: stock data output to json file
#
# The "股票" (stock) worksheet is missing a "property_category" column
# (always "stock") and a "legislator_id" column. Insert a new column
# before the existing "date" column, label it "property_category", and
# fill every data row with "stock". Also tidy up a few company names in
# column B that had a stray internal space before the trailing "司".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Column H currently holds "date" (name/owner/quantity/face_value/currency/
# total/date/legislator_name/legislator_id). Insert a fresh column there so
# date/legislator_name/legislator_id all shift one column to the right.
$ws.Columns.Item(8).Insert()

$ws.Cells.Item(1, 8).Value = "property_category"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}

# Remove the stray space in a handful of company names, e.g.
# "鴻海精密工業股份有限公 司" -> "鴻海精密工業股份有限公司"
$fixes = @{
    "鴻海精密工業股份有限公 司" = "鴻海精密工業股份有限公司"
    "碩禾電子材料股份有限公 司" = "碩禾電子材料股份有限公司"
    "欣陸投資控股股份有限公 司" = "欣陸投資控股股份有限公司"
    "興勤電子工業股份有限公 司" = "興勤電子工業股份有限公司"
}

for ($r = 2; $r -le 8; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null -and $fixes.ContainsKey($val)) {
        $cell.Value = $fixes[$val]
    }
}
